$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 holds a number-looking value but must stay text (like the original "—"),
# so force text entry with a leading apostrophe instead of letting Excel
# auto-detect it as a numeric value.
$ws.Range("A2").Value = "'1268"
$ws.Range("D2").Value = "31.439.951/0003-57"
$ws.Range("G2").Value = "582,00"
$ws.Range("H2").Value = "2025.05.30_PORTO SEGURO_211368.PDF"
